# feat: add 2022-Q3 data
#
# Inserts a new "2022-Q3" worksheet (right after "总计") populated with the
# quarter's fund-holdings table, and inserts a matching summary row at the
# top of the "总计" table so it lines up with the other quarters.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Add the new "2022-Q3" sheet right after the "总计" sheet.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)
$q3 = $wb.Worksheets.Add($null, $totalSheet)
$q3.Name = "2022-Q3"

# Header row (bold, bordered - matches the other quarter sheets' style).
$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"
$q3.Range("B1:H1").Font.Bold = $true
$q3.Range("B1:H1").HorizontalAlignment = -4108
$q3.Range("B1:H1").VerticalAlignment = -4160
$q3.Range("B1:H1").Borders.LineStyle = 1

# Fund holdings data for 2022-Q3.
# NB: each row literal is prefixed with the unary "," operator - without it
# PowerShell's array-subexpression flattens nested @(...) literals into one
# long list instead of a list of row-arrays.
$q3Rows = @(
    ,@(0,  "501208", "中欧创新未来混合（LOF）",                 "55.17", "85.30", "4.49",  "2.4771", 5)
    ,@(1,  "501081", "中欧科创主题混合（LOF）",                 "7.39",  "86.44", "8.51",  "0.6289", 3)
    ,@(2,  "015143", "中欧智能制造混合A",                       "0.96",  "84.58", "8.73",  "0.0838", 2)
    ,@(3,  "014837", "汇添富创新活力混合D",                     "2.46",  "70.85", "2.53",  "0.0622", 10)
    ,@(4,  "013369", "汇添富自主核心科技一年持有混合A",         "2.47",  "65.13", "2.48",  "0.0613", 10)
    ,@(5,  "002419", "汇添富创新活力混合A",                     "2.24",  "70.85", "2.53",  "0.0567", 10)
    ,@(6,  "015144", "中欧智能制造混合C",                       "0.51",  "84.58", "8.73",  "0.0445", 2)
    ,@(7,  "014737", "创金合信专精特新股票C",                   "0.58",  "81.21", "5.60",  "0.0325", 8)
    ,@(8,  "014736", "创金合信专精特新股票A",                   "0.35",  "81.21", "5.60",  "0.0196", 8)
    ,@(9,  "013370", "汇添富自主核心科技一年持有混合C",         "0.67",  "65.13", "2.48",  "0.0166", 10)
    ,@(10, "710002", "富安达策略精选混合",                     "0.59",  "50.67", "1.59",  "0.0094", 10)
    ,@(11, "014836", "汇添富创新活力混合C",                     "-0.01", "70.85", "2.53",  "-0.0003", 10)
)

$r = 2
foreach ($row in $q3Rows) {
    $q3.Cells.Item($r, 1).Value = $row[0]
    $q3.Cells.Item($r, 2).Value = "'" + $row[1]
    $q3.Cells.Item($r, 3).Value = $row[2]
    $q3.Cells.Item($r, 4).Value = "'" + $row[3]
    $q3.Cells.Item($r, 5).Value = "'" + $row[4]
    $q3.Cells.Item($r, 6).Value = "'" + $row[5]
    $q3.Cells.Item($r, 7).Value = "'" + $row[6]
    $q3.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2) Add a new summary row into "总计" for 2022-Q3, and push the existing
#    quarters (2022-Q2, 2022-Q1, 2021-Q4, 2021-Q3) down one row. Column A
#    is a plain 0-based row index, so it is left as-is for rows 2-5 and a
#    new "4" is written for the appended row 6 - only B/C/D shift.
#    (Range.Value *reads* aren't supported by this COM host, so the
#    "shifted" values are written as literals instead of copy-down.)
# ---------------------------------------------------------------------

# New row 6 - copy A5's number format/border over to A6 first.
$totalSheet.Range("A5").Copy()
$totalSheet.Range("A6").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$totalSheet.Range("A6").Value = 4
$totalSheet.Range("B6").Value = "2021-Q3"
$totalSheet.Range("C6").Value = 2
$totalSheet.Range("D6").Value = 0

$totalSheet.Range("B5").Value = "2021-Q4"
$totalSheet.Range("C5").Value = 1
$totalSheet.Range("D5").Value = 0.84

$totalSheet.Range("B4").Value = "2022-Q1"
$totalSheet.Range("C4").Value = 1
$totalSheet.Range("D4").Value = 0.6899999999999999

$totalSheet.Range("B3").Value = "2022-Q2"
$totalSheet.Range("C3").Value = 4
$totalSheet.Range("D3").Value = 2.61

$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 12
$totalSheet.Range("D2").Value = 3.49
